$d = $word.ActiveDocument

# The paragraph that currently reads "Update 1, 2, 3" (built from three
# separate runs: "Update 1", ", 2", ", 3", all Bold/size 22 half-points).
# We need to append a new run containing ", 4" with the same formatting.
$p = $d.Paragraphs(1)
$r = $p.Range

# $r.End currently points just past the trailing paragraph mark, so the
# real insertion point for new text is one character earlier.
$insertPos = $r.End - 1
$newText = ", 4"

# InsertAfter on the (un-collapsed) paragraph range correctly appends the
# text just before the paragraph mark, keeping paragraph boundaries intact.
$r.InsertAfter($newText)

# Grab a fresh Range that covers exactly the text we just inserted, then
# apply the same run formatting used by the rest of the line (bold, bold
# for complex scripts, 11pt / 22 half-points).
$newRange = $d.Range($insertPos, $insertPos + $newText.Length)
$newRange.Font.Bold = $true
$newRange.Font.BoldBi = $true
$newRange.Font.Size = 11
